$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.627927303314209
$ws.Range("B1").Value = 0.893427312374115
$ws.Range("C1").Value = 1.687890410423279
$ws.Range("D1").Value = 7.078995704650879
$ws.Range("E1").Value = 2.614516019821167
